# Insert a new weekly price record for "Ajo" (Garlic) at Terminal La Palmera
# de La Serena. The new observation is inserted as row 490, which pushes all
# subsequent rows (490-591) down by one (to 491-592), enlarging the used
# range from A1:R591 to A1:R592.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 490 (and everything below it) down by one row.
$ws.Rows("490:490").Insert()

# Populate the newly inserted row 490 with the new weekly record.
$ws.Range("A490").Value = 8
$ws.Range("B490").Value = "Terminal La Palmera de La Serena"
$ws.Range("C490").Value = "Coquimbo"
$ws.Range("D490").Value = 45275
$ws.Range("E490").Value = 4
$ws.Range("F490").Value = 100112003
$ws.Range("G490").Value = "Ajo"
$ws.Range("H490").Value = "Chino"
$ws.Range("I490").Value = "Primera"
$ws.Range("J490").Value = 400
$ws.Range("K490").Value = 23000
$ws.Range("L490").Value = 24000
$ws.Range("M490").Value = 23500
$ws.Range("N490").Value = "$/caja 10 kilos"
$ws.Range("O490").Value = "China"
$ws.Range("P490").Value = 2350
$ws.Range("Q490").Value = 10
$ws.Range("R490").Value = "Hortaliza"
